# Aerodynamic parameters.xlsx — add Gravity/Density/K.Viscosity/Total mass/
# Wind/ROC rows to the "Inputs" sheet (rows 10-15), match formatting of the
# existing parameter rows, give the kinematic-viscosity value a scientific
# number format, and move the active-cell selection to A16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs")

# Reference style (font + general look) used by the existing parameter rows
# (rows 4-9): a black-colored font, no fill. Reuse it for the new rows so the
# new cells look consistent with the rest of the table.
$refFont = $ws.Range("A9").Font.Color

function Set-ParamRow($row, $param, $symbol, $unit, $value, $formatValueCell) {
    $ws.Range("A" + $row).Value = $param
    $ws.Range("B" + $row).Value = $symbol
    $ws.Range("C" + $row).Value = $unit
    $ws.Range("D" + $row).Value = $value

    $ws.Range("A" + $row + ":C" + $row).Font.Color = $refFont
    if ($formatValueCell) {
        $ws.Range("D" + $row).Font.Color = $refFont
    }
}

Set-ParamRow 10 "Gravity"      "g0"     "m/s2"  9.80665    $true
Set-ParamRow 11 "Density"      "rho"    "kg/m3" 1.225      $true
Set-ParamRow 12 "K.Viscosity"  "mu"     "m2/s"  0.0000148  $false
Set-ParamRow 13 "Total mass"   "m_tot"  "kg"    7.08       $true
Set-ParamRow 14 "Wind"         "V_wind" "m/s"   3.4        $true
Set-ParamRow 15 "ROC"          "ROC"    "m/s"   2.8        $true

# Kinematic viscosity is tiny -> show it in scientific notation (D12 keeps
# the default font -- only the number format changes).
$ws.Range("D12").NumberFormat = "0.00E+00"

# Move the saved selection/active cell to A16 (matches the author's final
# cursor position after filling in the new rows).
$ws.Range("A16").Select()
